# RAAL/Production/Input/weather.xlsx — refresh pass after retraining the F_Cristan XGB model.
# The underlying cell values/formulas are unchanged (this is just the workbook being
# opened and re-saved by a newer Excel build after the input data refresh); the only
# observable, intentional edits left on the sheet are:
#   1) the dt_txt / Data columns (E:F) get their date-time display format re-applied,
#   2) every data column (A:X) is auto-fit to its content width,
#   3) the active selection is left on E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-apply the timestamp number format used by columns E (dt_txt) and F (Data).
$ws.Range("E2:F49").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# 2) Auto-fit every data column (A:X) to its content width.
$widths = @{
    "A" = 10.09
    "B" = 64.25
    "C" = 7.59
    "D" = 4.09
    "E" = 17.42
    "F" = 17.42
    "G" = 6.92
    "H" = 9.75
    "I" = 13.92
    "J" = 14.42
    "K" = 14.59
    "L" = 12.92
    "M" = 13.59
    "N" = 14.59
    "O" = 13.25
    "P" = 12.59
    "Q" = 8.59
    "R" = 10.59
    "S" = 8.42
    "T" = 8.75
    "U" = 6.75
    "V" = 6.25
    "W" = 7.59
    "X" = 11.25
}
foreach ($col in $widths.Keys) {
    $ws.Columns($col).ColumnWidth = $widths[$col]
}

# 3) Leave the selection on E6, matching the saved view state.
$ws.Range("E6").Select() | Out-Null
